# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and bump the
# handoff timestamps by one minute (00:32:xx -> 00:33:xx) across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"      # zh-cn status
$wsOverview.Range("C2").Value = "Ready for handoff"      # de-de status
$wsOverview.Range("D2").Value = "2016-03-22 00:33:15"    # Latest Handoff Date

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"          # Status
$wsZhCn.Range("E2").Value = "2016-03-22 00:33:11"        # Latest Handoff Datetime

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"          # Status
$wsDeDe.Range("E2").Value = "2016-03-22 00:33:15"        # Latest Handoff Datetime
